$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: phase plane -> trace-determinant plane
$ws.Range("D5").Value = "trace-determinant plane"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/17/trace_determinant_plane.html"

# Row 8: title update
$ws.Range("D8").Value = "제주어 기계번역 모델과 음성합성 모델에 관한 연구를 소개합니다."

# Row 9: title + link update
$ws.Range("D9").Value = "DS용 코딩에 대한 바른 이해 – Scientific programming이란?"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/scientific-programming-1/#utm_source=rss&utm_medium=rss&utm_campaign=scientific-programming-1"

# Row 28: title + link update
$ws.Range("D28").Value = '[4] Multi-Agent Actor-Critic for Mixed Cooperative-Competitive Environments (MADDPG)'
$ws.Range("E28").Value = "https://ropiens.tistory.com/124"

# Row 44: title + link update
$ws.Range("D44").Value = "5G 관련주 분석(1) - SK텔레콤"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/81"

# Row 51: title + link update
$ws.Range("D51").Value = "[MariaDB] 특정 열의 값이 NULL이 아닌 것만 보려면"
$ws.Range("E51").Value = "https://bskyvision.com/1197"
